# Add a "Link to Paper" hyperlink (OpenReview) on slide 2, and renumber the
# existing "Codebase" (GitHub) hyperlink relationship so the new link becomes
# rId2 and the GitHub link becomes rId3.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(2)
$sh = $s.Shapes.Item(2)

# Step 0: Clear the existing Codebase (GitHub) hyperlink address first so its
# relationship id (rId2) is freed up for the new OpenReview link.
$tr0 = $sh.TextFrame.TextRange
$codebaseLink0 = $tr0.Find("https://github.com/prob-ml/DynST?utm_source=chatgpt.com")
$asOld0 = $codebaseLink0.ActionSettings(1)
$asOld0.Hyperlink.Address = ""

# Step 1: Insert the new run's text right after the existing "Link to Paper: "
# run (appending to it keeps it as a single, unsplit run).
$tr1 = $sh.TextFrame.TextRange
$linkPara = $tr1.Find("Link to Paper: ")
$newRun = $linkPara.InsertAfter("https://openreview.net/forum?id=6quJeu5gJ7")

# Step 2: Assign the hyperlink to the newly inserted text first. This is what
# causes the text to split off into its own run (distinct from "Link to
# Paper: ") and freezes the paragraph's endParaRPr at the pre-edit (sz 1600)
# formatting.
$tr2 = $sh.TextFrame.TextRange
$sel = $tr2.Find("https://openreview.net/forum?id=6quJeu5gJ7")
$as = $sel.ActionSettings(1)
$as.Hyperlink.Address = "https://openreview.net/forum?id=6quJeu5gJ7"

# Step 3: Apply the run-level formatting to match the target: 18pt, not
# italic, single underline, blue (#1155CC), and an explicit (empty) effect
# list.
$tr3 = $sh.TextFrame.TextRange
$sel2 = $tr3.Find("https://openreview.net/forum?id=6quJeu5gJ7")
$sel2.Font.Size = 18
$sel2.Font.Italic = $false
$sel2.Font.Underline = $true
$sel2.Font.Color.RGB = 13391121
$sel2.Font.Shadow = $false

# Step 4: Re-apply the Codebase (GitHub) hyperlink address. Since rId2 is now
# used by the new OpenReview link, this picks up the next free id (rId3).
$tr4 = $sh.TextFrame.TextRange
$codebaseLink = $tr4.Find("https://github.com/prob-ml/DynST?utm_source=chatgpt.com")
$asOld = $codebaseLink.ActionSettings(1)
$asOld.Hyperlink.Address = "https://github.com/prob-ml/DynST?utm_source=chatgpt.com"
